# Update the "instructions" cell (A3) so the third line gains the trailing
# word "ลำดับ" after "7 ", matching the corrected task copy, and leave the
# active selection on A3 (matches the workbook's last-saved selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "`nในเกมนี้จะมีกล่องขึ้นมาทั้งหมด 9 กล่อง `nโดยจะมีกล่องกระพริบเป็นลำดับ`nจะเริ่มจาก 3 ลำดับไปจนถึง 7 ลำดับ`nแล้วต้องจำและคลิกตามลำดับนั้นๆ ให้ถูกต้อง`nแตะหน้าจอเพื่อดำเนินการต่อ"

$ws.Range("A3").Select()
